# A new weekly price record was added to the "Ajo" (garlic) price sheet.
# It is inserted as a new row 472, pushing the existing rows 472:508 down
# to 473:509 (dimension grows from A1:R508 to A1:R509).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 472, shifting rows 472:508
# (and all their data/styles) down to 473:509.
$ws.Rows.Item(472).Insert()

# Populate the new row 472 with the new record's data.
$ws.Range("A472").Value = 7
$ws.Range("B472").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C472").Value = "Ñuble"
$ws.Range("D472").Value = 45166
$ws.Range("E472").Value = 16
$ws.Range("F472").Value = 100112003
$ws.Range("G472").Value = "Ajo"
$ws.Range("H472").Value = "Chino"
$ws.Range("I472").Value = "Primera"
$ws.Range("J472").Value = 80
$ws.Range("K472").Value = 21000
$ws.Range("L472").Value = 21000
$ws.Range("M472").Value = 21000
$ws.Range("N472").Value = "$/caja 10 kilos"
$ws.Range("O472").Value = "China"
$ws.Range("P472").Value = 2100
$ws.Range("Q472").Value = 10
$ws.Range("R472").Value = "Hortaliza"
